$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 5436.2383
$ws.Range("I19").Value = 9715.546
$ws.Range("J19").Value = 729
$ws.Range("K19").Value = 9715.546
$ws.Range("L19").Value = 729
$ws.Range("M19").Value = -9540.546
$ws.Range("N19").Value = -1079
# Row 32
$ws.Range("H32").Value = 1148.6
$ws.Range("J32").Value = 1165.1111
$ws.Range("L32").Value = 1165.1111
$ws.Range("N32").Value = -1817.1111
# Row 41
$ws.Range("H41").Value = 294.64285
$ws.Range("I41").Value = 424
$ws.Range("J41").Value = 165.28572
$ws.Range("K41").Value = 424
$ws.Range("L41").Value = 165.28572
$ws.Range("M41").Value = 16
$ws.Range("N41").Value = -1045.28572
# Row 53
$ws.Range("H53").Value = 219.125
$ws.Range("I53").Value = 203.73334
$ws.Range("J53").Value = 450
$ws.Range("K53").Value = 203.73334
$ws.Range("L53").Value = 450
$ws.Range("M53").Value = 433.26666
$ws.Range("N53").Value = -1724
# Row 103
$ws.Range("H103").Value = 395.57144
$ws.Range("I103").Value = 236.5
$ws.Range("J103").Value = 459.2
$ws.Range("K103").Value = 709.5
$ws.Range("L103").Value = 1377.6
$ws.Range("M103").Value = -123.5
$ws.Range("N103").Value = -2549.6
# Row 129
$ws.Range("H129").Value = 38462290
$ws.Range("I129").Value = 125000300
$ws.Range("J129").Value = 947.7222
$ws.Range("K129").Value = 375000900
$ws.Range("L129").Value = 2843.1666
$ws.Range("M129").Value = -374995900
$ws.Range("N129").Value = -12843.1666
# Row 137
$ws.Range("H137").Value = 2250.8125
$ws.Range("I137").Value = 1906
$ws.Range("J137").Value = 2825.5
$ws.Range("K137").Value = 5718
$ws.Range("L137").Value = 8476.5
$ws.Range("M137").Value = -3168
$ws.Range("N137").Value = -13576.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 35
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2020.3226
$ws.Range("I20").Value = 2133.3809
$ws.Range("J20").Value = 1782.9
$ws.Range("K20").Value = 2133.3809
$ws.Range("L20").Value = 1782.9
$ws.Range("M20").Value = -1886.3809
$ws.Range("N20").Value = -2276.9
# Row 86
$ws.Range("H86").Value = 2088.7778
$ws.Range("I86").Value = 2072.9443
$ws.Range("J86").Value = 2120.4443
$ws.Range("K86").Value = 2072.9443
$ws.Range("L86").Value = 2120.4443
$ws.Range("M86").Value = -949.9443000000001
$ws.Range("N86").Value = -4366.4443
# Row 89
$ws.Range("H89").Value = 2088.7778
$ws.Range("I89").Value = 2072.9443
$ws.Range("J89").Value = 2120.4443
$ws.Range("K89").Value = 10364.7215
$ws.Range("L89").Value = 10602.2215
$ws.Range("M89").Value = -4748.7215
$ws.Range("N89").Value = -21834.2215
# Row 94
$ws.Range("H94").Value = 787.36
$ws.Range("I94").Value = 732
$ws.Range("J94").Value = 905
$ws.Range("K94").Value = 732
$ws.Range("L94").Value = 905
$ws.Range("M94").Value = -281
$ws.Range("N94").Value = -1807

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 1455432
$ws.Range("J68").Value = 1031.25
$ws.Range("L68").Value = 3093.75
$ws.Range("N68").Value = -4715.75
# Row 71
$ws.Range("H71").Value = 1455432
$ws.Range("J71").Value = 1031.25
$ws.Range("L71").Value = 9281.25
$ws.Range("N71").Value = -17393.25
# Row 131
$ws.Range("H131").Value = 1248.75
$ws.Range("J131").Value = 1367.3469
$ws.Range("L131").Value = 4102.0407
$ws.Range("N131").Value = -14182.0407
# Row 132
$ws.Range("H132").Value = 1093.8
$ws.Range("I132").Value = 872.9091
$ws.Range("K132").Value = 7856.1819
$ws.Range("M132").Value = -5326.1819
# Row 136
$ws.Range("H136").Value = 2561.25
$ws.Range("I136").Value = 2641.4285
$ws.Range("K136").Value = 7924.2855
$ws.Range("M136").Value = -2824.2855

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
# Row 83
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
# Row 113
$ws.Range("H113").Value = 2690.75
$ws.Range("I113").Value = 1916.6666
$ws.Range("K113").Value = 1916.6666
$ws.Range("M113").Value = 253.3334
# Row 132
$ws.Range("H132").Value = 3617.6667
$ws.Range("I132").Value = 3803
$ws.Range("J132").Value = 3580.6
$ws.Range("K132").Value = 11409
$ws.Range("L132").Value = 10741.8
$ws.Range("M132").Value = -8879
$ws.Range("N132").Value = -15801.8

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 2414.8572
$ws.Range("I16").Value = 1200
$ws.Range("K16").Value = 1200
$ws.Range("M16").Value = -1030
# Row 46
$ws.Range("H46").Value = 38206.297
$ws.Range("I46").Value = 72283
$ws.Range("J46").Value = 1508.3077
$ws.Range("K46").Value = 72283
$ws.Range("L46").Value = 1508.3077
$ws.Range("M46").Value = -72095
$ws.Range("N46").Value = -1884.3077
# Row 55
$ws.Range("H55").Value = 230.9375
$ws.Range("I55").Value = 260.66666
$ws.Range("J55").Value = 141.75
$ws.Range("K55").Value = 260.66666
$ws.Range("L55").Value = 141.75
$ws.Range("M55").Value = -87.66665999999998
$ws.Range("N55").Value = -487.75
# Row 80
$ws.Range("H80").Value = 17441.6
$ws.Range("J80").Value = 17441.6
$ws.Range("L80").Value = 17441.6
$ws.Range("N80").Value = -19687.6
# Row 82
$ws.Range("H82").Value = 1984.8334
$ws.Range("I82").Value = 1354
$ws.Range("J82").Value = 2489.5
$ws.Range("K82").Value = 1354
$ws.Range("L82").Value = 2489.5
$ws.Range("M82").Value = -993
$ws.Range("N82").Value = -3211.5
# Row 83
$ws.Range("H83").Value = 17441.6
$ws.Range("J83").Value = 17441.6
$ws.Range("L83").Value = 52324.8
$ws.Range("N83").Value = -63556.8
# Row 85
$ws.Range("H85").Value = 1984.8334
$ws.Range("I85").Value = 1354
$ws.Range("J85").Value = 2489.5
$ws.Range("K85").Value = 1354
$ws.Range("L85").Value = 2489.5
$ws.Range("M85").Value = -106
$ws.Range("N85").Value = -4985.5
# Row 132
$ws.Range("H132").Value = 7629.4897
$ws.Range("I132").Value = 11251.038
$ws.Range("J132").Value = 3535.5652
$ws.Range("K132").Value = 33753.114
$ws.Range("L132").Value = 10606.6956
$ws.Range("M132").Value = -31223.114
$ws.Range("N132").Value = -15666.6956

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 3942.3818
$ws.Range("I136").Value = 594.7778
$ws.Range("J136").Value = 5570.946
$ws.Range("K136").Value = 1784.3334
$ws.Range("L136").Value = 16712.838
$ws.Range("M136").Value = 765.6666
$ws.Range("N136").Value = -21812.838
